$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 3 (shifts rows 3-36 down to 4-37,
# which already reproduces the required shift of every subsequent row,
# including the style change that lands on the new row 31 and the new
# trailing row 37).
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new pre-requisite question.
$ws.Range("B3").Value2 = "Is it possible to get a service endpoint created in the SANDBOX to create the end to end devops pipeline"
$ws.Rows.Item(3).RowHeight = 30

# The numbering that used to live in A3:A6 (1,2,3,4) is removed - those
# cells become blank (matching the rest of the A column) while keeping
# their existing style.
$ws.Range("A4").ClearContents()
$ws.Range("A5").ClearContents()
$ws.Range("A6").ClearContents()
$ws.Range("A7").ClearContents()

# Update the active selection to reflect the edited cell.
$ws.Range("B4").Select() | Out-Null
